$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as text so values like "62.80" or
# "37.401.94" are not reinterpreted as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.401.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.037.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.33%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.647"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.80"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +15.93%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +8.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.63"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0754"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.83%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.908"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +9.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.333.39"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.56"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +8.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.31"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +21.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.035.14"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.181.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.28"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0872"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.33"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.78"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +27.33%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.52"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.80"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.114"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +31.31%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.33%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +9.45%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +10.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.66"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +11.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0614"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +16.57%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.02"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +27.58%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.106"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +23.01%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.30%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0219"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.80%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.14"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.24%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.10"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +13.26%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.68"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +24.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.93"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +13.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.07"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +8.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.432.21"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.95"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.62"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.58%  "
